$wb = $excel.ActiveWorkbook

# Rename the two sheets.
$ws1 = $wb.Worksheets.Item("룸배정표(8.25~26) (2)")
$ws1.Name = "룸배정표 1차수 (8.25~26) (2)"

$ws2 = $wb.Worksheets.Item("룸배정표 (8.27~28)")
$ws2.Name = "룸배정표 2차수  (8.27~28)"

# Update the window placement for the workbook.
$excel.Left = 3384
$excel.Top = 3360
$excel.Width = 17280
$excel.Height = 8880

# Adjust the first sheet's view: clear topLeftCell scroll position and change the selection.
$ws1.Activate()
$ws1.Range("A1:E43").Select()
